# ASP.NET MVC Working with Data.pptx - apply commit "Added homeworks as README's"
#
# Summary of the change:
#  - The "Homework" slide (slide 43, id 464) is removed from the deck (its
#    homework content moved to a README instead of living on a slide).
#  - What used to be slide 44 ("Free Trainings @ Telerik Academy") shifts up
#    to become slide 43, and its cached slide-number field text is updated.
#  - A handful of unrelated small text tidy-ups: runs that were split across
#    multiple <a:r> for no reason get merged back into a single run on a few
#    slides, and the cached "today" date fields in the handout/notes masters
#    get refreshed.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Remove the "Homework" slide (slide 43).
# ---------------------------------------------------------------------
$p.Slides.Item(43).Delete()

# ---------------------------------------------------------------------
# 2. The slide that used to be #44 is now #43 - refresh its cached slide
#    number field text so it reads "43" instead of "44".
# ---------------------------------------------------------------------
$sLast = $p.Slides.Item($p.Slides.Count)
for ($i = 1; $i -le $sLast.Shapes.Count; $i++) {
    $shp = $sLast.Shapes.Item($i)
    if ($shp.Name -like "Slide Number Placeholder*") {
        $shp.TextFrame.TextRange.Text = [string]$sLast.SlideIndex
    }
}

# ---------------------------------------------------------------------
# 3. Merge runs that were needlessly split across multiple <a:r> elements.
# ---------------------------------------------------------------------

# Slide 2 - "Scaffolding" slide: "Data " + "Validation" -> "Data Validation"
$s2 = $p.Slides.Item(2)
foreach ($shp in $s2.Shapes) {
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -like "*Data Validation*") {
            $shp.TextFrame.TextRange.Text = $shp.TextFrame.TextRange.Text
        }
    }
}

# Slide 5 - "Demo: Create Scaffold": merge "Create " / "CRUD pages " / "with
# read/write actions, using Entity Framework" into a single run.
$s5 = $p.Slides.Item(5)
foreach ($shp in $s5.Shapes) {
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -like "Create CRUD pages *") {
            $shp.TextFrame.TextRange.Text = "Create CRUD pages with read/write actions, using Entity Framework"
        }
    }
}

# Slide 27 - jQuery validation bullet: merge the two runs back together.
$s27 = $p.Slides.Item(27)
foreach ($shp in $s27.Shapes) {
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -like "jQuery validation library*") {
            $shp.TextFrame.TextRange.Text = "jQuery validation library required for unobtrusive JavaScript validation"
        }
    }
}

Write-Output "done"
